$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the entire existing data row (row 2), which held a sample
# "Nivel 0" dataset, and replace it with a single remaining value
# (AH2), matching the new reduced template contents.
$ws.Rows.Item(2).Clear()

$ws.Range("AH2").Value = "Maria Teresa del Socorro Cazola Bravo"
